# SFS_AllMobile_Tbl_Live.xlsx update
# "All softwares updated working fine" - refresh the appointment/date/voice/form
# fields in row 2 to the latest (26/12/2023) run, replacing the previous
# (30/11/2023) run's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# O2 -> AP_Value (appointment date/time range message)
$ws.Range("O2").Value = "date_range`nAppointment Date : 26/12/2023, Time : [ 09:10 AM to 09:14 AM ]"

# AC2 -> DT_Value (selected date)
$ws.Range("AC2").Value = "26/12/2023"

# AS2 -> VR_Value (voice recording filename)
$ws.Range("AS2").Value = "voice_record_26122023"

# AV2 -> UF_Prefix (uploaded form filename)
$ws.Range("AV2").Value = "formshow_26122023"
